$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 3
$ws.Range("F3").Value = 15
$ws.Range("H3").Value = 18

# Row 5
$ws.Range("F5").Value = 15
$ws.Range("H5").Value = 19

# Row 8
$ws.Range("E8").Value = 46

# Row 11
$ws.Range("E11").Value = 20
$ws.Range("F11").Value = 13
$ws.Range("H11").Value = 14

# Row 15
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 49
$ws.Range("H15").Value = 60

# Row 16
$ws.Range("F16").Value = 93
$ws.Range("H16").Value = 180

# Row 18
$ws.Range("E18").Value = 98
$ws.Range("F18").Value = 32
$ws.Range("H18").Value = 55
